$d = $word.ActiveDocument

$d.Content.Find.Execute("HC254", $true, $false, $false, $false, $false,
                         $true, 1, $false, "BC344", 2)
